$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.903.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.897.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.7557"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'240.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3044"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.08%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'25.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.48%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.06834"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.22%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07987"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.7489"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.35%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.904.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.211"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.85%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'91.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'29.906.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'13.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.85%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'5.949"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'240.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.000007727"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.83%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'2.155.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'6.954"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.38%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'9.241"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.48%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'165.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'18.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.1299"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.71%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.018"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.73%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.434"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.16%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'4.282"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'4.027"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.48%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.05373"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.36%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.256"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.72%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.7258"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.37%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.723"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.01920"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.786"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'6.171"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.30%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.4405"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.11%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'72.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.93%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.909"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.89%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.8272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'101.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'7.571"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.18%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'9.796"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'2.054.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'36.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.81%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.05966"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.58%  '
$ws.Range("E51").Style = "Normal"
